# Update "想去人数" (interested-count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet, matching the regenerated data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 190
$ws1.Range("F3").Value = 5373
$ws1.Range("F7").Value = 608
$ws1.Range("F8").Value = 581
$ws1.Range("F9").Value = 1052
$ws1.Range("F11").Value = 1473
$ws1.Range("F12").Value = 4303
$ws1.Range("F13").Value = 440
$ws1.Range("F14").Value = 194
$ws1.Range("F17").Value = 3442
$ws1.Range("F19").Value = 1095
$ws1.Range("F20").Value = 105
$ws1.Range("F21").Value = 42
$ws1.Range("F22").Value = 201
$ws1.Range("F26").Value = 72
$ws1.Range("F27").Value = 312
$ws1.Range("F28").Value = 32
$ws1.Range("F32").Value = 27

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 190
$ws4.Range("F4").Value = 5373
$ws4.Range("F8").Value = 608
$ws4.Range("F9").Value = 581
$ws4.Range("F10").Value = 1052
$ws4.Range("F12").Value = 1473
$ws4.Range("F13").Value = 4304
$ws4.Range("F14").Value = 440
$ws4.Range("F15").Value = 194
$ws4.Range("F18").Value = 3442
$ws4.Range("F20").Value = 1095
$ws4.Range("F21").Value = 105
$ws4.Range("F22").Value = 42
$ws4.Range("F23").Value = 201
$ws4.Range("F27").Value = 72
$ws4.Range("F28").Value = 312
$ws4.Range("F29").Value = 32
$ws4.Range("F33").Value = 27
